# Update "Yearly" sheet: December row (row 14) values for 401K (E) and Suzie's Roth IRA (F).
# Totals (row 15) and the dependent cells on the "All Time" sheet are formulas and will
# recalculate automatically.
$wb = $excel.ActiveWorkbook

$wsYearly = $wb.Worksheets.Item("Yearly")
$wsYearly.Range("E14").Value = 40.29
$wsYearly.Range("F14").Value = 49.22

# Restore the selection on the "Yearly" sheet.
$wsYearly.Range("I14").Select()

# "All Time" sheet: clear the frozen/scrolled top-left cell and move the selection.
$wsAllTime = $wb.Worksheets.Item("All Time")
$wsAllTime.Activate()
$wsAllTime.Application.ActiveWindow.ScrollRow = 1
$wsAllTime.Application.ActiveWindow.ScrollColumn = 1
$wsAllTime.Range("K12").Select()
